$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Footer 1 (Word OM index) -> docPr id="2" (PearsonLogo) ---
# rename inline picture "image2.png" -> "image1.png"
$f1 = $sec.Footers.Item(1)
$s1 = $f1.Range.InlineShapes.Item(1)
$shp1 = $s1.ConvertToShape()
$shp1.Name = "image1.png"
$null = $shp1.ConvertToInlineShape()

# --- Footer 2 (Word OM index) -> docPr id="3" (PearsonLogo) ---
# rename inline picture "image2.png" -> "image1.png"
$f2 = $sec.Footers.Item(2)
$s2 = $f2.Range.InlineShapes.Item(1)
$shp2 = $s2.ConvertToShape()
$shp2.Name = "image1.png"
$null = $shp2.ConvertToInlineShape()

# --- Header 2 (Word OM index) -> docPr id="1" (BTec_Logo-Orange) ---
# rename inline picture "image1.jpg" -> "image2.jpg"
$h2 = $sec.Headers.Item(2)
$s3 = $h2.Range.InlineShapes.Item(1)
$shp3 = $s3.ConvertToShape()
$shp3.Name = "image2.jpg"
$null = $shp3.ConvertToInlineShape()

Write-Host "Renamed inline shapes in footers and header."
